# Trade #96 closed at 2026-02-16 21:39:07 - momentum DOWN +0.000%
#
# This edit reflects two trading-log events:
#   1) An existing OPEN "leadlag" trade (#60, row 49 on the "leadlag" sheet /
#      row 61 on "All Trades") gets closed out: exit price, P&L, exit reason
#      and duration are filled in, and status flips to CLOSED. The Summary
#      and Comparison roll-up sheets are updated to match the new totals.
#   2) A brand-new "momentum" trade (#96) is opened and appended as a new
#      row at the bottom of the "momentum" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value as literal text (never let Excel's automatic
# type-inference turn a date-looking / percent-looking string into a
# real date serial or percentage number).
# ---------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Drop the temporary "@" text format again so the cell is left with the
    # default (unstyled) formatting, matching the rest of the sheet.
    $cell.ClearFormats()
}

# =======================================================================
# 1) Summary sheet — OVERALL + leadlag roll-up rows
# =======================================================================
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Cells.Item(2, 3).Value = 60            # C2 Total Trades (overall)
Set-TextCell $wsSummary 2 4 "68.3%"               # D2 Win Rate
Set-TextCell $wsSummary 2 5 "+16.3031%"           # E2 Total P&L %
Set-TextCell $wsSummary 2 6 "+0.2717%"            # F2 Avg Trade

$wsSummary.Cells.Item(3, 3).Value = 71            # C3 Total Trades (leadlag)
Set-TextCell $wsSummary 3 4 "43.7%"               # D3 Win Rate
Set-TextCell $wsSummary 3 5 "+11.4523%"           # E3 Total P&L %
Set-TextCell $wsSummary 3 6 "+0.1613%"            # F3 Avg Trade

# =======================================================================
# 2) leadlag sheet — row 49 (trade #60) closes out
# =======================================================================
$wsLeadlag = $wb.Worksheets.Item("leadlag")

$wsLeadlag.Cells.Item(49, 7).Value = 68398.489224   # G49 Exit Price
Set-TextCell $wsLeadlag 49 8 "CLOSED"               # H49 Status
$wsLeadlag.Cells.Item(49, 9).Value = 0.5401         # I49 P&L %
$wsLeadlag.Cells.Item(49, 10).Value = 5.4           # J49 P&L $
Set-TextCell $wsLeadlag 49 13 "time_exit_5min"      # M49 Exit Reason
$wsLeadlag.Cells.Item(49, 14).Value = 5             # N49 Duration (min)

# =======================================================================
# 3) All Trades sheet — append the now-closed trade #60 as new row 61
# =======================================================================
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Cells.Item(61, 1).Value = 60
Set-TextCell $wsAll 61 2 "2026-02-16"
Set-TextCell $wsAll 61 3 "21:34:02"
Set-TextCell $wsAll 61 4 "leadlag"
Set-TextCell $wsAll 61 5 "DOWN"
$wsAll.Cells.Item(61, 6).Value = 68769.89999999999
$wsAll.Cells.Item(61, 7).Value = 68398.489224
Set-TextCell $wsAll 61 8 "CLOSED"
$wsAll.Cells.Item(61, 9).Value = 0.5401
$wsAll.Cells.Item(61, 10).Value = 5.4
$wsAll.Cells.Item(61, 11).Value = 0.75
Set-TextCell $wsAll 61 12 "Coinbase leading with -0.086% move"
Set-TextCell $wsAll 61 13 "time_exit_5min"
$wsAll.Cells.Item(61, 14).Value = 5

# =======================================================================
# 4) momentum sheet — append newly-opened trade #96 as new row 26
# =======================================================================
$wsMomentum = $wb.Worksheets.Item("momentum")

$wsMomentum.Cells.Item(26, 1).Value = 96
Set-TextCell $wsMomentum 26 2 "2026-02-16"
Set-TextCell $wsMomentum 26 3 "21:39:07"
Set-TextCell $wsMomentum 26 4 "momentum"
Set-TextCell $wsMomentum 26 5 "DOWN"
$wsMomentum.Cells.Item(26, 6).Value = 68298.875
Set-TextCell $wsMomentum 26 8 "OPEN"
$wsMomentum.Cells.Item(26, 9).Value = 0
$wsMomentum.Cells.Item(26, 10).Value = 0
$wsMomentum.Cells.Item(26, 11).Value = 0.9
Set-TextCell $wsMomentum 26 12 "Downward momentum: -0.242% over 10 samples"
$wsMomentum.Cells.Item(26, 14).Value = 0

# =======================================================================
# 5) Comparison sheet — leadlag roll-up row
# =======================================================================
$wsComparison = $wb.Worksheets.Item("Comparison")

$wsComparison.Cells.Item(2, 2).Value = 71          # B2 Total Trades
Set-TextCell $wsComparison 2 3 "43.7%"             # C2 Win Rate
Set-TextCell $wsComparison 2 4 "3.16"              # D2 Profit Factor
